# Realestate Update resale numbers 2023-06-30 11:36
# Appends a new data row (row 94) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 94

# Columns A-D hold text values (Date/Time/Weekday/Week are stored as text
# in this sheet, not as native date/number types). Prefix with a leading
# apostrophe so Excel stores them as text instead of auto-converting the
# date/time/number-looking strings into a date serial or numeric value.
$ws.Cells.Item($row, 1).Value = "'2023-06-30"
$ws.Cells.Item($row, 2).Value = "'11:36:55"
$ws.Cells.Item($row, 3).Value = "'Friday"
$ws.Cells.Item($row, 4).Value = "'26"

# Columns E-T hold the numeric resale counts for each city.
$ws.Cells.Item($row, 5).Value  = 123386   # Beijing
$ws.Cells.Item($row, 6).Value  = 134499   # Guangzhou
$ws.Cells.Item($row, 7).Value  = 163489   # Suzhou
$ws.Cells.Item($row, 8).Value  = 134071   # Hangzhou
$ws.Cells.Item($row, 9).Value  = 177030   # Nanjing
$ws.Cells.Item($row, 10).Value = 115484   # Xi_an
$ws.Cells.Item($row, 11).Value = 204719   # Chengdu
$ws.Cells.Item($row, 12).Value = 226062   # Chongqing
$ws.Cells.Item($row, 13).Value = 176634   # Tianjin
$ws.Cells.Item($row, 14).Value = 104652   # Hefei
$ws.Cells.Item($row, 15).Value = 39850    # Fuzhou
$ws.Cells.Item($row, 16).Value = 33705    # Xiamen
$ws.Cells.Item($row, 17).Value = 52632    # Changsha
$ws.Cells.Item($row, 18).Value = -1       # Shanghai
$ws.Cells.Item($row, 19).Value = 35916    # Shenzhen
$ws.Cells.Item($row, 20).Value = -1       # Wuhan
